$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B133").Value = "Test annualized geometric excess return"
$ws.Range("C133").Value = "return_annualized_excess_test1"
$ws.Range("A133").Value = "return_annualized_excess1"

$ws.Range("A134").Value = "return_annualized_excess2"
$ws.Range("B134").Value = "Test annualized arithmetic excess return"
$ws.Range("C134").Value = "return_annualized_excess_test2"
